# Edit applied per commit "Fixed update to excel issue":
# - Add a new "PO Forecast" sheet (sheetId=3) with forecast data
#   (columns: ds, PO_Forecast, yhat_lower, yhat_upper)
# - Rename "Requested quantity" header to "Weekly_PO_Qty" on the "Weekly Quantity" sheet
# - Rename "Requested quantity" header to "Monthly_PO_Qty" on the "Monthly Trend" sheet

$wb = $excel.ActiveWorkbook

# --- Rename header on "Weekly Quantity" sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Rename header on "Monthly Trend" sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end of the workbook ---
# Duplicate "Weekly Quantity" (instead of Worksheets.Add) so the new sheet
# inherits the same sheet-level properties (outline/page-setup prefs, page
# margins) and cell formatting (bold/centered/bordered header, date number
# format on the first data column), then wipe the copied values.
$wsWeekly.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"
$wsForecast.Cells.ClearContents()

# The source header only spans 2 columns (A:B) and the source date column
# only spans 27 rows; extend the copied formatting to cover the new sheet's
# 4 columns x 35 rows.
$wsForecast.Range("B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)

$wsForecast.Range("A27").Copy()
$wsForecast.Range("A28:A35").PasteSpecial(-4122)

# --- Header values ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Forecast data rows ---
$forecastData = @(
    @(44955.99999999999, 146, 1.742700319204755, 299.5027582159902),
    @(44983.99999999999, 142, -7.91495752978489, 296.618642445348),
    @(44990.99999999999, 140, -4.065772799681369, 294.2293153867589),
    @(44997.99999999999, 139, -4.672937277533527, 290.4571462186574),
    @(45053.99999999999, 130, -17.98733595226494, 286.6170471390983),
    @(45060.99999999999, 129, -25.61433297902525, 270.5968345971293),
    @(45081.99999999999, 125, -20.90087505150275, 279.4175676135313),
    @(45088.99999999999, 124, -27.11822640731911, 278.1040966604496),
    @(45095.99999999999, 123, -25.10622276487808, 276.1053561270215),
    @(45102.99999999999, 122, -29.95034554401544, 264.0067428668272),
    @(45109.99999999999, 120, -33.37080927800945, 279.7738199878333),
    @(45116.99999999999, 119, -17.762926105863, 276.0911012050643),
    @(45123.99999999999, 118, -17.21126533357818, 260.481079571683),
    @(45130.99999999999, 117, -24.87987306796238, 266.6310572163071),
    @(45137.99999999999, 116, -29.39029888329701, 268.6724471800283),
    @(45144.99999999999, 115, -37.64241942215545, 259.8824633545863),
    @(45151.99999999999, 113, -41.02230942711368, 246.832135401425),
    @(45158.99999999999, 112, -41.01058670995335, 262.5233692013923),
    @(45165.99999999999, 111, -34.1200820377242, 261.1108385403078),
    @(45172.99999999999, 110, -45.23092224099088, 251.001138087466),
    @(45186.99999999999, 108, -45.07958232931986, 254.8537946224166),
    @(45193.99999999999, 106, -41.13390247322449, 264.6785961497245),
    @(45200.99999999999, 105, -56.88234085808905, 256.1356372521067),
    @(45207.99999999999, 104, -59.42409418475212, 267.2873858973067),
    @(45214.99999999999, 103, -39.79082804388727, 252.938207325325),
    @(45221.99999999999, 102, -46.99545307342268, 247.6984871665373),
    @(45228.99999999999, 101, -47.96564523212126, 244.9198230429763),
    @(45235.99999999999, 99, -49.62666654554188, 241.1668252003537),
    @(45242.99999999999, 98, -42.529474837824, 248.4386608509636),
    @(45249.99999999999, 97, -54.30830112626904, 240.3339152246517),
    @(45256.99999999999, 96, -54.11204482030026, 248.721827269527),
    @(45263.99999999999, 95, -48.02457881948936, 234.495749184773),
    @(45270.99999999999, 93, -37.37844070387779, 243.3475237270203),
    @(45277.99999999999, 92, -50.84908313537846, 253.7631750006227)
)

$r = 2
foreach ($row in $forecastData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}
